# Scheduled data refresh: update crypto price/volume snapshot (cryptos.xlsx)
# generated by the "Updated cryptos list ... with GitHub Actions" workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$addr, [string]$val) {
    # The source cells are plain text (prices use "." as a thousands
    # separator, e.g. "44.142.20", volumes are padded percentage
    # strings). Force text format so numeric-looking values (e.g.
    # "318.39") are not reinterpreted as floating point numbers, then
    # restore the default "Normal" style so no stray number format
    # sticks to the cell (matches the original unstyled cells).
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '44.142.20'
Set-TextValue 'E2' '  +1.71%  '
Set-TextValue 'D3' '2.247.81'
Set-TextValue 'E3' '  +0.57%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '318.39'
Set-TextValue 'D6' '101.13'
Set-TextValue 'E6' '  +1.74%  '
Set-TextValue 'D7' '0.574'
Set-TextValue 'E7' '  -1.43%  '
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'D9' '0.544'
Set-TextValue 'E9' '  -3.28%  '
Set-TextValue 'D10' '36.86'
Set-TextValue 'E10' '  -1.11%  '
Set-TextValue 'D11' '0.0827'
Set-TextValue 'E11' '  -0.28%  '
Set-TextValue 'E12' '  -2.54%  '
Set-TextValue 'E13' '  -2.00%  '
Set-TextValue 'D14' '2.589.74'
Set-TextValue 'E14' '  +0.49%  '
Set-TextValue 'D15' '2.281.32'
Set-TextValue 'E15' '  +1.50%  '
Set-TextValue 'D16' '0.849'
Set-TextValue 'E16' '  -2.00%  '
Set-TextValue 'D17' '14.18'
Set-TextValue 'E17' '  -1.43%  '
Set-TextValue 'D18' '43.998.60'
Set-TextValue 'E18' '  +1.48%  '
Set-TextValue 'D19' '13.46'
Set-TextValue 'E19' '  -5.15%  '
Set-TextValue 'D20' '0.0₃0977'
Set-TextValue 'E20' '  +0.55%  '
Set-TextValue 'E21' '  -3.04%  '
Set-TextValue 'D22' '65.53'
Set-TextValue 'E22' '  +0.40%  '
Set-TextValue 'D23' '3.10'
Set-TextValue 'E23' '  -3.80%  '
Set-TextValue 'D24' '234.83'
Set-TextValue 'E24' '  -0.61%  '
Set-TextValue 'D25' '2.06'
Set-TextValue 'E25' '  -5.81%  '
Set-TextValue 'E26' '  -0.08%  '
Set-TextValue 'D27' '10.45'
Set-TextValue 'E27' '  +3.35%  '
Set-TextValue 'D28' '2.22'
Set-TextValue 'E28' '  +0.05%  '
Set-TextValue 'E29' '  +1.98%  '
Set-TextValue 'D30' '6.11'
Set-TextValue 'E30' '  -4.62%  '
Set-TextValue 'B31' 'Monero'
Set-TextValue 'C31' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D31' '158.39'
Set-TextValue 'E31' '  +0.39%  '
Set-TextValue 'B32' 'EthereumClassic'
Set-TextValue 'C32' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D32' '20.07'
Set-TextValue 'E32' '  -1.21%  '
Set-TextValue 'E33' '  -3.47%  '
Set-TextValue 'D34' '2.67'
Set-TextValue 'E34' '  -1.34%  '
Set-TextValue 'D35' '3.17'
Set-TextValue 'E35' '  -2.05%  '
Set-TextValue 'B36' 'ARBITRUM'
Set-TextValue 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D36' '1.95'
Set-TextValue 'E36' '  +3.01%  '
Set-TextValue 'B37' 'Kaspa'
Set-TextValue 'C37' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D37' '0.112'
Set-TextValue 'E37' '  +7.18%  '
Set-TextValue 'E38' '  -2.18%  '
Set-TextValue 'D39' '16.08'
Set-TextValue 'E39' '  +10.87%  '
Set-TextValue 'D40' '3.68'
Set-TextValue 'E40' '  -0.62%  '
Set-TextValue 'D41' '4.15'
Set-TextValue 'E41' '  -6.10%  '
Set-TextValue 'D42' '0.0315'
Set-TextValue 'E42' '  -2.33%  '
Set-TextValue 'E43' '  +0.05%  '
Set-TextValue 'D44' '1.747.36'
Set-TextValue 'E44' '  -3.40%  '
Set-TextValue 'D45' '0.197'
Set-TextValue 'E45' '  -3.18%  '
Set-TextValue 'B46' 'BitcoinSV'
Set-TextValue 'C46' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D46' '82.23'
Set-TextValue 'E46' '  -2.73%  '
Set-TextValue 'B47' 'ordi'
Set-TextValue 'C47' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 'D47' '74.68'
Set-TextValue 'E47' '  +0.37%  '
Set-TextValue 'D48' '5.15'
Set-TextValue 'B49' 'Stacks'
Set-TextValue 'C49' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D49' '1.68'
Set-TextValue 'E49' '  +3.93%  '
Set-TextValue 'B50' 'Aave'
Set-TextValue 'C50' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D50' '102.72'
Set-TextValue 'E50' '  -0.94%  '
Set-TextValue 'D51' '57.41'
Set-TextValue 'E51' '  -2.24%  '
